$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5169141.5
$ws.Range("I17").Value = 3766.6667
$ws.Range("J17").Value = 6137649.5
$ws.Range("K17").Value = 11300.0001
$ws.Range("L17").Value = 18412948.5
$ws.Range("M17").Value = -11132.0001
$ws.Range("N17").Value = -18413284.5
$ws.Range("H19").Value = 1442.25
$ws.Range("I19").Value = 1442.25
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1442.25
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -1267.25
$ws.Range("N19").ClearContents()
$ws.Range("H33").Value = 1694
$ws.Range("I33").Value = 738.82355
$ws.Range("K33").Value = 738.82355
$ws.Range("M33").Value = -509.82355
$ws.Range("H100").Value = 52273.668
$ws.Range("I100").Value = 81407.38
$ws.Range("J100").Value = 4931.375
$ws.Range("K100").Value = 81407.38
$ws.Range("L100").Value = 4931.375
$ws.Range("M100").Value = -80866.38
$ws.Range("N100").Value = -6013.375
$ws.Range("H113").Value = 7845.6665
$ws.Range("I113").Value = 7185
$ws.Range("J113").Value = 7977.8
$ws.Range("K113").Value = 7185
$ws.Range("L113").Value = 7977.8
$ws.Range("M113").Value = -3931
$ws.Range("N113").Value = -14485.8
$ws.Range("H132").Value = 1932.7727
$ws.Range("I132").Value = 1932.7727
$ws.Range("K132").Value = 5798.3181
$ws.Range("M132").Value = -3268.3181
$ws.Range("H137").Value = 8298.771000000001
$ws.Range("I137").Value = 3243
$ws.Range("J137").Value = 17651.95
$ws.Range("K137").Value = 9729
$ws.Range("L137").Value = 52955.85000000001
$ws.Range("M137").Value = -7179
$ws.Range("N137").Value = -58055.85000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4301.237
$ws.Range("I2").Value = 3325.0293
$ws.Range("K2").Value = 3325.0293
$ws.Range("M2").Value = -3212.0293
$ws.Range("H32").Value = 3138.7292
$ws.Range("I32").Value = 1625.9269
$ws.Range("J32").Value = 11999.429
$ws.Range("K32").Value = 1625.9269
$ws.Range("L32").Value = 11999.429
$ws.Range("M32").Value = -1338.9269
$ws.Range("N32").Value = -12573.429
$ws.Range("H46").Value = 20261.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 20261.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 20261.5
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -20899.5
$ws.Range("H61").Value = 5330.6562
$ws.Range("I61").Value = 3641.1155
$ws.Range("J61").Value = 12652
$ws.Range("K61").Value = 3641.1155
$ws.Range("L61").Value = 12652
$ws.Range("M61").Value = -3429.1155
$ws.Range("N61").Value = -13076
$ws.Range("H116").Value = 4301.237
$ws.Range("I116").Value = 3325.0293
$ws.Range("K116").Value = 3325.0293
$ws.Range("M116").Value = -1031.0293
$ws.Range("H122").Value = 2211.0715
$ws.Range("I122").Value = 2175.4167
$ws.Range("J122").Value = 2425
$ws.Range("K122").Value = 6526.250100000001
$ws.Range("L122").Value = 7275
$ws.Range("M122").Value = -4076.250100000001
$ws.Range("N122").Value = -12175
$ws.Range("H132").Value = 3475788.2
$ws.Range("I132").Value = 3707333.2
$ws.Range("K132").Value = 11121999.6
$ws.Range("M132").Value = -11119469.6
$ws.Range("H136").Value = 5330.6562
$ws.Range("I136").Value = 3641.1155
$ws.Range("J136").Value = 12652
$ws.Range("K136").Value = 10923.3465
$ws.Range("L136").Value = 37956
$ws.Range("M136").Value = -8373.3465
$ws.Range("N136").Value = -43056

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4301.237
$ws.Range("I3").Value = 3325.0293
$ws.Range("K3").Value = 3325.0293
$ws.Range("M3").Value = -3211.0293
$ws.Range("H80").Value = 292
$ws.Range("I80").Value = 292
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 292
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 706
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 292
$ws.Range("I83").Value = 292
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 1460
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 3532
$ws.Range("N83").ClearContents()
$ws.Range("H105").Value = 2869.1914
$ws.Range("I105").Value = 3607.9473
$ws.Range("J105").Value = 2367.8928
$ws.Range("K105").Value = 3607.9473
$ws.Range("L105").Value = 2367.8928
$ws.Range("M105").Value = -1860.9473
$ws.Range("N105").Value = -5861.8928
$ws.Range("H134").Value = 10128173
$ws.Range("I134").Value = 37050370
$ws.Range("J134").Value = 32349
$ws.Range("K134").Value = 111151110
$ws.Range("L134").Value = 97047
$ws.Range("M134").Value = -111148575
$ws.Range("N134").Value = -102117

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6545.7856
$ws.Range("I58").Value = 4149.364
$ws.Range("J58").Value = 15332.667
$ws.Range("K58").Value = 4149.364
$ws.Range("L58").Value = 15332.667
$ws.Range("M58").Value = -3946.364
$ws.Range("N58").Value = -15738.667
$ws.Range("H86").Value = 5222.222
$ws.Range("I86").Value = 5333.3335
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 5333.3335
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -4210.3335
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 5222.222
$ws.Range("I89").Value = 5333.3335
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 26666.6675
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -21050.6675
$ws.Range("N89").Value = -36232
$ws.Range("H94").Value = 1047.8636
$ws.Range("I94").Value = 1386.2
$ws.Range("J94").Value = 948.35297
$ws.Range("K94").Value = 1386.2
$ws.Range("L94").Value = 948.35297
$ws.Range("M94").Value = -935.2
$ws.Range("N94").Value = -1850.35297
$ws.Range("H99").Value = 10186.8
$ws.Range("I99").Value = 6467.769
$ws.Range("J99").Value = 11977.444
$ws.Range("K99").Value = 6467.769
$ws.Range("L99").Value = 11977.444
$ws.Range("M99").Value = -4969.769
$ws.Range("N99").Value = -14973.444
$ws.Range("H122").Value = 2643.4167
$ws.Range("I122").Value = 2378
$ws.Range("K122").Value = 7134
$ws.Range("M122").Value = -4684
$ws.Range("H126").Value = 10186.8
$ws.Range("I126").Value = 6467.769
$ws.Range("J126").Value = 11977.444
$ws.Range("K126").Value = 19403.307
$ws.Range("L126").Value = 35932.33199999999
$ws.Range("M126").Value = -16933.307
$ws.Range("N126").Value = -40872.33199999999
$ws.Range("H132").Value = 28704.04
$ws.Range("I132").Value = 18930.312
$ws.Range("J132").Value = 46079.555
$ws.Range("K132").Value = 56790.936
$ws.Range("L132").Value = 138238.665
$ws.Range("M132").Value = -54260.936
$ws.Range("N132").Value = -143298.665
$ws.Range("H134").Value = 21082.416
$ws.Range("I134").Value = 19299.1
$ws.Range("K134").Value = 57897.3
$ws.Range("M134").Value = -55362.3
$ws.Range("H136").Value = 6545.7856
$ws.Range("I136").Value = 4149.364
$ws.Range("J136").Value = 15332.667
$ws.Range("K136").Value = 12448.092
$ws.Range("L136").Value = 45998.001
$ws.Range("M136").Value = -9898.091999999999
$ws.Range("N136").Value = -51098.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1788.3334
$ws.Range("I134").Value = 1788.3334
$ws.Range("K134").Value = 5365.0002
$ws.Range("M134").Value = -295.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7701.4194
$ws.Range("I70").Value = 7248.8184
$ws.Range("K70").Value = 7248.8184
$ws.Range("M70").Value = -6978.8184
$ws.Range("H73").Value = 7701.4194
$ws.Range("I73").Value = 7248.8184
$ws.Range("K73").Value = 7248.8184
$ws.Range("M73").Value = -6312.8184
$ws.Range("H102").Value = 2514.7896
$ws.Range("I102").Value = 2538.3572
$ws.Range("J102").Value = 2448.8
$ws.Range("K102").Value = 2538.3572
$ws.Range("L102").Value = 2448.8
$ws.Range("M102").Value = -916.3571999999999
$ws.Range("N102").Value = -5692.8
$ws.Range("H113").Value = 401999.6
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 1999998
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 1999998
$ws.Range("M113").Value = -330
$ws.Range("N113").Value = -2004338
$ws.Range("H126").Value = 3113.4736
$ws.Range("I126").Value = 2867.6
$ws.Range("J126").Value = 4035.5
$ws.Range("K126").Value = 8602.799999999999
$ws.Range("L126").Value = 12106.5
$ws.Range("M126").Value = -6132.799999999999
$ws.Range("N126").Value = -17046.5
$ws.Range("H132").Value = 12339.286
$ws.Range("I132").Value = 13986.723
$ws.Range("J132").Value = 2454.6667
$ws.Range("K132").Value = 41960.169
$ws.Range("L132").Value = 7364.000100000001
$ws.Range("M132").Value = -39430.169
$ws.Range("N132").Value = -12424.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 39999
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H46").Value = 2713.2856
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 2999
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 2999
$ws.Range("M46").Value = -811
$ws.Range("N46").Value = -3375
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20676
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22340
$ws.Range("H122").Value = 3796.5
$ws.Range("I122").Value = 3701.3333
$ws.Range("K122").Value = 11103.9999
$ws.Range("M122").Value = -8653.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10325.429
$ws.Range("I122").Value = 6098
$ws.Range("J122").Value = 13496
$ws.Range("K122").Value = 18294
$ws.Range("L122").Value = 40488
$ws.Range("M122").Value = -15844
$ws.Range("N122").Value = -45388
$ws.Range("H136").Value = 744.9524
$ws.Range("I136").Value = 769.2895
$ws.Range("J136").Value = 513.75
$ws.Range("K136").Value = 2307.8685
$ws.Range("L136").Value = 1541.25
$ws.Range("M136").Value = 242.1315
$ws.Range("N136").Value = -6641.25
